$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update vm_pu values for rows 2-25, columns B-F and I-N
# (case with 380 kV done)

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.041573871308939
$ws.Cells.Item(2, 4).Value = 1.042700715924044
$ws.Cells.Item(2, 5).Value = 1.050016769059546
$ws.Cells.Item(2, 6).Value = 1.061391677084075
$ws.Cells.Item(2, 9).Value = 1.04126628743465
$ws.Cells.Item(2, 10).Value = 1.046654347908918
$ws.Cells.Item(2, 11).Value = 1.045476661899094
$ws.Cells.Item(2, 12).Value = 1.05277221746129
$ws.Cells.Item(2, 13).Value = 1.064115863913117
$ws.Cells.Item(2, 14).Value = 1.019421011382735

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.042521533354295
$ws.Cells.Item(3, 4).Value = 1.04339750614697
$ws.Cells.Item(3, 5).Value = 1.050890560040428
$ws.Cells.Item(3, 6).Value = 1.062419276010029
$ws.Cells.Item(3, 9).Value = 1.041507185954621
$ws.Cells.Item(3, 10).Value = 1.04724810323241
$ws.Cells.Item(3, 11).Value = 1.0459846359876
$ws.Cells.Item(3, 12).Value = 1.053458212026311
$ws.Cells.Item(3, 13).Value = 1.064957537986113
$ws.Cells.Item(3, 14).Value = 1.019621419241437

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.043135150904895
$ws.Cells.Item(4, 4).Value = 1.043848696072581
$ws.Cells.Item(4, 5).Value = 1.051456731556323
$ws.Cells.Item(4, 6).Value = 1.063085252187491
$ws.Cells.Item(4, 9).Value = 1.041662112161525
$ws.Cells.Item(4, 10).Value = 1.047632092363693
$ws.Cells.Item(4, 11).Value = 1.0463129729291
$ws.Cells.Item(4, 12).Value = 1.0539022260982
$ws.Cells.Item(4, 13).Value = 1.06550259705826
$ws.Cells.Item(4, 14).Value = 1.019750945779433

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.043393214614813
$ws.Cells.Item(5, 4).Value = 1.044038451663088
$ws.Cells.Item(5, 5).Value = 1.05169493326036
$ws.Cells.Item(5, 6).Value = 1.063365479007642
$ws.Cells.Item(5, 9).Value = 1.04172701490775
$ws.Cells.Item(5, 10).Value = 1.047793470174912
$ws.Cells.Item(5, 11).Value = 1.046450919612542
$ws.Cells.Item(5, 12).Value = 1.054088919698227
$ws.Cells.Item(5, 13).Value = 1.065731844094879
$ws.Cells.Item(5, 14).Value = 1.019805362273262

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.043436550420239
$ws.Cells.Item(6, 4).Value = 1.044070316842892
$ws.Cells.Item(6, 5).Value = 1.051734939082109
$ws.Cells.Item(6, 6).Value = 1.063412544975624
$ws.Cells.Item(6, 9).Value = 1.041737898966316
$ws.Cells.Item(6, 10).Value = 1.047820563190398
$ws.Cells.Item(6, 11).Value = 1.04647407638936
$ws.Cells.Item(6, 12).Value = 1.054120268111159
$ws.Cells.Item(6, 13).Value = 1.065770341760992
$ws.Cells.Item(6, 14).Value = 1.019814496892975

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.043138598777355
$ws.Cells.Item(7, 4).Value = 1.043851231301613
$ws.Cells.Item(7, 5).Value = 1.051459913700826
$ws.Cells.Item(7, 6).Value = 1.063088995609391
$ws.Cells.Item(7, 9).Value = 1.04166298029261
$ws.Cells.Item(7, 10).Value = 1.047634248904261
$ws.Cells.Item(7, 11).Value = 1.046314816519108
$ws.Cells.Item(7, 12).Value = 1.053904720590221
$ws.Cells.Item(7, 13).Value = 1.065505659860252
$ws.Cells.Item(7, 14).Value = 1.019751673039027

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.041894051282091
$ws.Cells.Item(8, 4).Value = 1.042936132404471
$ws.Cells.Item(8, 5).Value = 1.05031191057938
$ws.Cells.Item(8, 6).Value = 1.061738740641816
$ws.Cells.Item(8, 9).Value = 1.041347896876574
$ws.Cells.Item(8, 10).Value = 1.046855053306987
$ws.Cells.Item(8, 11).Value = 1.045648407359087
$ws.Cells.Item(8, 12).Value = 1.053004025093537
$ws.Cells.Item(8, 13).Value = 1.064400219936295
$ws.Cells.Item(8, 14).Value = 1.019488771145569

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.039704231846259
$ws.Cells.Item(9, 4).Value = 1.041326118846772
$ws.Cells.Item(9, 5).Value = 1.048294935949314
$ws.Cells.Item(9, 6).Value = 1.05936751951944
$ws.Cells.Item(9, 9).Value = 1.040785416435291
$ws.Cells.Item(9, 10).Value = 1.045480439640843
$ws.Cells.Item(9, 11).Value = 1.044471419471694
$ws.Cells.Item(9, 12).Value = 1.05141792088551
$ws.Cells.Item(9, 13).Value = 1.062455707453114
$ws.Cells.Item(9, 14).Value = 1.019024366350745

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.038246577683845
$ws.Cells.Item(10, 4).Value = 1.040254537861833
$ws.Cells.Item(10, 5).Value = 1.046954360715351
$ws.Cells.Item(10, 6).Value = 1.05779222027015
$ws.Cells.Item(10, 9).Value = 1.04040557564446
$ws.Cells.Item(10, 10).Value = 1.044563025686858
$ws.Cells.Item(10, 11).Value = 1.043685004027953
$ws.Cells.Item(10, 12).Value = 1.05036127185251
$ws.Cells.Item(10, 13).Value = 1.06116172297164
$ws.Cells.Item(10, 14).Value = 1.018714019955692

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.037615935762436
$ws.Cells.Item(11, 4).Value = 1.039790964961879
$ws.Cells.Item(11, 5).Value = 1.046374858109119
$ws.Cells.Item(11, 6).Value = 1.057111420418269
$ws.Cells.Item(11, 9).Value = 1.040239955214969
$ws.Cells.Item(11, 10).Value = 1.044165547636561
$ws.Cells.Item(11, 11).Value = 1.043344071439333
$ws.Cells.Item(11, 12).Value = 1.049903921700962
$ws.Cells.Item(11, 13).Value = 1.060601985227196
$ws.Cells.Item(11, 14).Value = 1.018579464817077

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.037381767987164
$ws.Cells.Item(12, 4).Value = 1.039618839004599
$ws.Cells.Item(12, 5).Value = 1.046159752735959
$ws.Cells.Item(12, 6).Value = 1.056858739716452
$ws.Cells.Item(12, 9).Value = 1.040178264570405
$ws.Cells.Item(12, 10).Value = 1.04401787259398
$ws.Cells.Item(12, 11).Value = 1.043217373200194
$ws.Cells.Item(12, 12).Value = 1.049734070326208
$ws.Cells.Item(12, 13).Value = 1.060394159739103
$ws.Cells.Item(12, 14).Value = 1.018529459395529

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.03743199410459
$ws.Cells.Item(13, 4).Value = 1.039655757615563
$ws.Cells.Item(13, 5).Value = 1.046205886866333
$ws.Cells.Item(13, 6).Value = 1.056912931576628
$ws.Cells.Item(13, 9).Value = 1.040191505191914
$ws.Cells.Item(13, 10).Value = 1.044049550917131
$ws.Cells.Item(13, 11).Value = 1.043244553135862
$ws.Cells.Item(13, 12).Value = 1.04977050270314
$ws.Cells.Item(13, 13).Value = 1.060438735102755
$ws.Cells.Item(13, 14).Value = 1.018540186883311

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.037596577729934
$ws.Cells.Item(14, 4).Value = 1.039776735629943
$ws.Cells.Item(14, 5).Value = 1.046357074412138
$ws.Cells.Item(14, 6).Value = 1.057090529693361
$ws.Cells.Item(14, 9).Value = 1.040234859349915
$ws.Cells.Item(14, 10).Value = 1.044153341460241
$ws.Cells.Item(14, 11).Value = 1.043333599756736
$ws.Cells.Item(14, 12).Value = 1.049889881139824
$ws.Cells.Item(14, 13).Value = 1.060584804546763
$ws.Cells.Item(14, 14).Value = 1.018575331874574

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.037697993852659
$ws.Cells.Item(15, 4).Value = 1.039851282914997
$ws.Cells.Item(15, 5).Value = 1.046450245655176
$ws.Cells.Item(15, 6).Value = 1.057199980113837
$ws.Cells.Item(15, 9).Value = 1.040261548520499
$ws.Cells.Item(15, 10).Value = 1.044217285759059
$ws.Cells.Item(15, 11).Value = 1.04338845630082
$ws.Cells.Item(15, 12).Value = 1.049963437982926
$ws.Cells.Item(15, 13).Value = 1.060674814181828
$ws.Cells.Item(15, 14).Value = 1.018596982475093

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.038288442586434
$ws.Cells.Item(16, 4).Value = 1.040285312807782
$ws.Cells.Item(16, 5).Value = 1.046992841043132
$ws.Cells.Item(16, 6).Value = 1.057837430515483
$ws.Cells.Item(16, 9).Value = 1.040416543203079
$ws.Cells.Item(16, 10).Value = 1.044589400183261
$ws.Cells.Item(16, 11).Value = 1.043707622043789
$ws.Cells.Item(16, 12).Value = 1.050391628669378
$ws.Cells.Item(16, 13).Value = 1.061198882932925
$ws.Cells.Item(16, 14).Value = 1.018722946320304

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.038658958454751
$ws.Cells.Item(17, 4).Value = 1.040557683957168
$ws.Cells.Item(17, 5).Value = 1.04733345873787
$ws.Cells.Item(17, 6).Value = 1.058237639176823
$ws.Cells.Item(17, 9).Value = 1.040513460549354
$ws.Cells.Item(17, 10).Value = 1.044822756262426
$ws.Cells.Item(17, 11).Value = 1.043907717228758
$ws.Cells.Item(17, 12).Value = 1.050660271775982
$ws.Cells.Item(17, 13).Value = 1.061527769664147
$ws.Cells.Item(17, 14).Value = 1.018801914018523

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.038875125364012
$ws.Cells.Item(18, 4).Value = 1.040716594723949
$ws.Cells.Item(18, 5).Value = 1.04753222930196
$ws.Cells.Item(18, 6).Value = 1.058471200934243
$ws.Cells.Item(18, 9).Value = 1.040569880111321
$ws.Cells.Item(18, 10).Value = 1.044958846499193
$ws.Cells.Item(18, 11).Value = 1.044024389797175
$ws.Cells.Item(18, 12).Value = 1.050816984625712
$ws.Cells.Item(18, 13).Value = 1.061719658281839
$ws.Cells.Item(18, 14).Value = 1.018847957815288

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.038948841366339
$ws.Cells.Item(19, 4).Value = 1.040770786153304
$ws.Cells.Item(19, 5).Value = 1.047600020840516
$ws.Cells.Item(19, 6).Value = 1.058550860934762
$ws.Cells.Item(19, 9).Value = 1.040589098935649
$ws.Cells.Item(19, 10).Value = 1.045005245925922
$ws.Cells.Item(19, 11).Value = 1.044064165391367
$ws.Cells.Item(19, 12).Value = 1.0508704226611
$ws.Cells.Item(19, 13).Value = 1.061785096564594
$ws.Cells.Item(19, 14).Value = 1.018863654712374

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.038619200299537
$ws.Cells.Item(20, 4).Value = 1.040528456831378
$ws.Cells.Item(20, 5).Value = 1.047296903948844
$ws.Cells.Item(20, 6).Value = 1.058194687449864
$ws.Cells.Item(20, 9).Value = 1.040503073680799
$ws.Cells.Item(20, 10).Value = 1.044797721676242
$ws.Cells.Item(20, 11).Value = 1.043886252995755
$ws.Cells.Item(20, 12).Value = 1.050631447057987
$ws.Cells.Item(20, 13).Value = 1.061492477580751
$ws.Cells.Item(20, 14).Value = 1.018793443259284

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.03754810971991
$ws.Cells.Item(21, 4).Value = 1.039741108797515
$ws.Cells.Item(21, 5).Value = 1.046312549361604
$ws.Cells.Item(21, 6).Value = 1.057038225981268
$ws.Cells.Item(21, 9).Value = 1.040222097370953
$ws.Cells.Item(21, 10).Value = 1.044122778664081
$ws.Cells.Item(21, 11).Value = 1.043307379415151
$ws.Cells.Item(21, 12).Value = 1.049854726362318
$ws.Cells.Item(21, 13).Value = 1.060541788353331
$ws.Cells.Item(21, 14).Value = 1.018564983255885

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.036875139942831
$ws.Cells.Item(22, 4).Value = 1.039246452334087
$ws.Cells.Item(22, 5).Value = 1.045694501673035
$ws.Cells.Item(22, 6).Value = 1.056312262777855
$ws.Cells.Item(22, 9).Value = 1.040044442413296
$ws.Cells.Item(22, 10).Value = 1.043698218916301
$ws.Cells.Item(22, 11).Value = 1.042943067864527
$ws.Cells.Item(22, 12).Value = 1.049366538517445
$ws.Cells.Item(22, 13).Value = 1.059944550551938
$ws.Cells.Item(22, 14).Value = 1.018421193145324

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.03723184938558
$ws.Cells.Item(23, 4).Value = 1.039508642602257
$ws.Cells.Item(23, 5).Value = 1.046022059001919
$ws.Cells.Item(23, 6).Value = 1.05669700023239
$ws.Cells.Item(23, 9).Value = 1.040138714733062
$ws.Cells.Item(23, 10).Value = 1.043923304402391
$ws.Cells.Item(23, 11).Value = 1.043136229280847
$ws.Cells.Item(23, 12).Value = 1.049625319968732
$ws.Cells.Item(23, 13).Value = 1.060261110052058
$ws.Cells.Item(23, 14).Value = 1.018497432936418

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.038637165140119
$ws.Cells.Item(24, 4).Value = 1.040541663184026
$ws.Cells.Item(24, 5).Value = 1.047313421194105
$ws.Cells.Item(24, 6).Value = 1.058214095094681
$ws.Cells.Item(24, 9).Value = 1.040507767401638
$ws.Cells.Item(24, 10).Value = 1.044809033797463
$ws.Cells.Item(24, 11).Value = 1.043895951881061
$ws.Cells.Item(24, 12).Value = 1.050644471652159
$ws.Cells.Item(24, 13).Value = 1.061508424385624
$ws.Cells.Item(24, 14).Value = 1.018797270882476

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.040269964253344
$ws.Cells.Item(25, 4).Value = 1.041742041006325
$ws.Cells.Item(25, 5).Value = 1.048815659331894
$ws.Cells.Item(25, 6).Value = 1.059979570014796
$ws.Cells.Item(25, 9).Value = 1.040931688804734
$ws.Cells.Item(25, 10).Value = 1.045835991128318
$ws.Cells.Item(25, 11).Value = 1.044776012494752
$ws.Cells.Item(25, 12).Value = 1.051827838166098
$ws.Cells.Item(25, 13).Value = 1.06295799961386
$ws.Cells.Item(25, 14).Value = 1.019144558725134

